# LOQ4261.docx edit script
# Applies a content rotation among several paragraphs plus a structural
# reshuffle that moves the "Bibliografia" heading (and its list paragraph)
# from before the bibliography text to after it, and appends a new
# paragraph with the instructor's name at the very end of that block.

$d = $word.ActiveDocument

function Replace-InParagraph($ParaIndex, $OldText, $NewText) {
    $p = $d.Paragraphs.Item($ParaIndex)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $ok = $r.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        throw "Replace-InParagraph failed at paragraph $ParaIndex for text: $OldText"
    }
}

# ---------------------------------------------------------------------
# Step 1: simple text rotations (paragraph positions/styles unchanged)
# ---------------------------------------------------------------------

# "Objetivos" (PT) paragraph gets the old "Programa resumido" (PT) text
Replace-InParagraph 6 `
    "Apresentar um quadro conceitual de análise para auxiliar na formulação, avaliação e desenvolvimento de modelos para Planejamento, Programação e Controle da Produção nos diferentes ambientes de produção." `
    "Caracterização do planejamento e controle da produção. 2. Planejamento agregado da produção. 3. Planejamento mestre da produção. 4. Planejamento e controle de estoques. 5. Planejamento de recursos de materiais (MRP). 6. Programação detalhada da produção. 7. Sistema MRPII e Sistema ERP. 8.Tambor-Pulmão-Corda - OPT. 9. Teoria das Restrições (TOC)."

# "Objetivos" (EN) paragraph gets the old "Programa resumido" (EN) text
Replace-InParagraph 7 `
    "To present a conceptual framework of analysis to assist in the formulation, evaluation and development of models for Planning, Programming and Production Control in different production environments." `
    "Characterization of production programming and control. 2. Aggregate Production Planning. 3. Master Production Schedulling. 4. Inventory planning and control. 5. Material Requirement Planning (MRP). 6. Detailed scheduling of production. 7. Production control systems. 8. Drum-Buffer-Rope – Opt; 9. Theory of Constraints (TOC)"

# "Docente(s)" list paragraph gets the old "Objetivos" (PT) text
Replace-InParagraph 9 `
    "8971158 - Claudemir Leif Tramarico" `
    "Apresentar um quadro conceitual de análise para auxiliar na formulação, avaliação e desenvolvimento de modelos para Planejamento, Programação e Controle da Produção nos diferentes ambientes de produção."

# "Programa resumido" (EN) paragraph gets the old "Objetivos" (EN) text
Replace-InParagraph 12 `
    "Characterization of production programming and control. 2. Aggregate Production Planning. 3. Master Production Schedulling. 4. Inventory planning and control. 5. Material Requirement Planning (MRP). 6. Detailed scheduling of production. 7. Production control systems. 8. Drum-Buffer-Rope – Opt; 9. Theory of Constraints (TOC)" `
    "To present a conceptual framework of analysis to assist in the formulation, evaluation and development of models for Planning, Programming and Production Control in different production environments."

# "Programa" (PT) paragraph gets the old "Método" value text
Replace-InParagraph 14 `
    "Caracterização do planejamento e controle da produção. 2. Planejamento agregado da produção. 3. Planejamento mestre da produção. 4. Planejamento e controle de estoques. 5. Planejamento de recursos de materiais (MRP). 6. Programação detalhada da produção. 7. Sistema MRPII e Sistema ERP. 8.Tambor-Pulmão-Corda - OPT. 9. Teoria das Restrições (TOC)." `
    "Provas, atividades em grupo e atividades individuais."

# ---------------------------------------------------------------------
# Step 2: rotate the three labeled values inside the "Avaliação" block
# (paragraph 17) -- Método's value <- old Critério value,
# Critério's value <- old Norma de recuperação value. Each value is
# located by anchoring right after its bold label and up to the next
# line break, so the (temporarily) duplicated text introduced by the
# rotation can never create an ambiguous match.
# ---------------------------------------------------------------------

function Get-ValueRangeAfterLabel($ParaIndex, $LabelText) {
    $p = $d.Paragraphs.Item($ParaIndex)
    $scan = $d.Range($p.Range.Start, $p.Range.End)
    $ok = $scan.Find.Execute($LabelText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Get-ValueRangeAfterLabel: label not found: $LabelText"
    }
    $valStart = $scan.End
    $p2 = $d.Paragraphs.Item($ParaIndex)
    $scan2 = $d.Range($valStart, $p2.Range.End)
    $okVtab = $scan2.Find.Execute("`v", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($okVtab) {
        $valEnd = $scan2.Start
    } else {
        $p3 = $d.Paragraphs.Item($ParaIndex)
        $valEnd = $p3.Range.End - 1
    }
    return $d.Range($valStart, $valEnd)
}

$metodoValue = Get-ValueRangeAfterLabel 17 "Método: "
$metodoValue.Text = "Média das atividades avaliativas"

$criterioValue = Get-ValueRangeAfterLabel 17 "Critério: "
$criterioValue.Text = "MF = (0,5 M + 0,5 R) M = Média de aproveitamento do aluno, antes da recuperação R = Nota de uma prova de recuperação MF = nota final de aproveitamento, após a recuperação Aprovação com média final de aproveitamento maior ou igual a 5,0. A recuperação deverá consistir em uma prova escrita englobando a matéria toda do semestre. Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%."

# ---------------------------------------------------------------------
# Step 3: remove the old "MF = ..." run (now orphaned, directly after
# the "Norma de recuperação: " label) then splice in the bibliography
# list text (formerly its own "Normal" paragraph) right there, keeping
# the "Avaliação" paragraph's ListBullet style throughout.
# ---------------------------------------------------------------------

# Paragraph 18 = "Bibliografia" heading; paragraph 19 = the bibliography
# list text (Normal style). Give the list paragraph the ListBullet style
# so that, once it's merged backwards, the combined paragraph keeps the
# ListBullet formatting that paragraph 17 already uses.
$biblioListPara = $d.Paragraphs.Item(19)
$biblioListPara.Style = "ListBullet"

# Delete the old "MF = ..." run's text (the one directly following
# "Norma de recuperação: " inside paragraph 17).
$p17 = $d.Paragraphs.Item(17)
$scan = $d.Range($p17.Range.Start, $p17.Range.End)
$okLabel = $scan.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okLabel) {
    throw "Could not find 'Norma de recuperação: ' label"
}
$oldMfStart = $scan.End
$p17After = $d.Paragraphs.Item(17)
$oldMfEnd = $p17After.Range.End - 1   # exclude the paragraph mark
$mfRange = $d.Range($oldMfStart, $oldMfEnd)
$mfRange.Delete()

# Now delete the "Bibliografia" heading paragraph (text + its own
# paragraph mark); this merges the (now ListBullet) bibliography-list
# paragraph up into its place.
$bibHeading = $d.Paragraphs.Item(18)
$bibHeadingRange = $d.Range($bibHeading.Range.Start, $bibHeading.Range.End)
$bibHeadingRange.Delete()

# Finally, merge that bibliography-list paragraph into paragraph 17 by
# deleting paragraph 17's own paragraph mark (the combined paragraph
# keeps paragraph 18's -- i.e. the list's -- ListBullet style, since
# that is the style we just applied to it).
$p17Final = $d.Paragraphs.Item(17)
$markRange = $d.Range($p17Final.Range.End - 1, $p17Final.Range.End)
$markRange.Delete()

# ---------------------------------------------------------------------
# Step 4: insert a fresh "Bibliografia" heading paragraph and a
# paragraph with the instructor reference, right after the merged
# paragraph (and before "Requisitos").
# ---------------------------------------------------------------------

# "Requisitos" is now paragraph 18 (right after the merged paragraph
# that absorbed the bibliography list). Insert the two new paragraphs
# in front of it -- they initially inherit its Heading2 style, so fix
# up the second one's style afterwards.
$reqPara = $d.Paragraphs.Item(18)
$reqPara.Range.InsertBefore("Bibliografia`r8971158 - Claudemir Leif Tramarico`r")

$headingPara = $d.Paragraphs.Item(18)
$headingPara.Style = "Heading2"

$namePara = $d.Paragraphs.Item(19)
$namePara.Style = "Normal"
